$d = $word.ActiveDocument

function Get-ParaContaining($doc, $needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $t = $paras.Item($i).Range.Text
        if ($t.Contains($needle)) {
            return $paras.Item($i)
        }
    }
    return $null
}

# -----------------------------------------------------------------
# Change 1: append ", Excel" to the Business Tools line, as its own run.
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Vim, IntelliJ, Eclipse, Android Studio, JIRA, Rally", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(", Excel")

# -----------------------------------------------------------------
# Change 2: "Master's in Business Administration" -> "Master" + " in Business
# Administration" (i.e. the "'s" is dropped). Word's automatic "_GoBack"
# bookmark - which marks the location of the most recent edit - moves along
# with it: from the end of the document (after the Honors line) to right
# after the word "Master", where this edit happens.
# -----------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Master’s in Business Administration", $true, $false, $false, $false, $false, $true, 1, $false, "Master in Business Administration", 2) | Out-Null

# Locate the paragraph that now reads "Master in Business Administration, May
# 2019" and compute the character offset right after the word "Master".
$eduPara = Get-ParaContaining $d "Master in Business Administration"
$pRange = $eduPara.Range
$pText = $pRange.Text
$afterMaster = $pRange.Start + ($pText.IndexOf("Master") + "Master".Length)

# Re-seat (and thereby relocate) the "_GoBack" bookmark right between "Master"
# and " in Business Administration". Adding a bookmark with an existing name
# replaces the old one, so this automatically removes it from its old spot at
# the end of the Honors paragraph.
$bmPoint = $d.Range($afterMaster, $afterMaster)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# The find/replace above collapsed "Master in Business Administration, May
# 2019" into a single run (aside from the split the bookmark insertion just
# introduced). Re-establish the original, separate run boundaries for ", "
# and "May 2019" that followed "Master's in Business Administration" before
# this edit.
$eduPara2 = Get-ParaContaining $d "in Business Administration"
$pRange2 = $eduPara2.Range
$pText2 = $pRange2.Text
$tail = ", May 2019"
$tailIdx = $pText2.IndexOf($tail)
$tailRange = $d.Range($pRange2.Start + $tailIdx, $pRange2.Start + $tailIdx + $tail.Length)
$tailRange.Text = ""

$eduPara3 = Get-ParaContaining $d "in Business Administration"
$pRange3 = $eduPara3.Range
$endPoint = $d.Range($pRange3.End - 1, $pRange3.End - 1)
$endPoint.InsertAfter(", ")

$eduPara4 = Get-ParaContaining $d "in Business Administration"
$pRange4 = $eduPara4.Range
$endPoint2 = $d.Range($pRange4.End - 1, $pRange4.End - 1)
$endPoint2.InsertAfter("May 2019")
